$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204, shifting the existing rows 204..324 down to 205..325.
$ws.Rows.Item(204).Insert()

# Populate the newly-inserted row 204 with the latest weekly price entry
# for "Ajo / Chino / Primera" at "Vega Monumental Concepción".
$ws.Range("A204").Value = 11
$ws.Range("B204").Value = "Vega Monumental Concepción"
$ws.Range("C204").Value = "Bíobío"
$ws.Range("D204").Value = 45176
$ws.Range("E204").Value = 8
$ws.Range("F204").Value = 100112003
$ws.Range("G204").Value = "Ajo"
$ws.Range("H204").Value = "Chino"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 300
$ws.Range("K204").Value = 18000
$ws.Range("L204").Value = 19000
$ws.Range("M204").Value = 18667
$ws.Range("N204").Value = "$/caja 10 kilos"
$ws.Range("O204").Value = "China"
$ws.Range("P204").Value = 1867
$ws.Range("Q204").Value = 10
$ws.Range("R204").Value = "Hortaliza"
